$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.534982800483704
$ws.Range("B1").Value = 2.337340593338013
$ws.Range("C1").Value = 2.823202610015869
$ws.Range("D1").Value = 3.279439687728882
$ws.Range("E1").Value = 2.169921398162842
